$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update fixture text (results now known) and mark pick outcome in column G
$ws.Range("A2").Value = "Cape Verde ✓ - Eswatini: 3:0"
$ws.Range("G2").Value = "✓"

$ws.Range("A3").Value = "Iceland - France : 2:2"

$ws.Range("A4").Value = "Cameroon  - Angola: 0:0"

$ws.Range("A5").Value = "North Macedonia  - Kazakhstan: 1:1"

$ws.Range("A6").Value = "Sweden X - Kosovo: 0:1"
$ws.Range("G6").Value = "X"

$ws.Range("A7").Value = "Tunisia ✓ - Namibia: 3:0"
$ws.Range("G7").Value = "✓"

# Remove the old row 8 fixture entirely (shrinks used range to A1:G7)
$ws.Range("A8:G8").Delete()
